$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 72; existing rows 72..135 shift down to 73..136.
$ws.Rows.Item(72).EntireRow.Insert()

# Populate the newly inserted row 72 with its data.
$r = 72
$ws.Cells.Item($r, 1).Value = 1
$ws.Cells.Item($r, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item($r, 3).Value = "Arica y Parinacota"
$ws.Cells.Item($r, 4).Value = 44467
$ws.Cells.Item($r, 5).Value = 15
$ws.Cells.Item($r, 6).Value = "Fruta"
$ws.Cells.Item($r, 7).Value = 100108
$ws.Cells.Item($r, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item($r, 9).Value = 100108006
$ws.Cells.Item($r, 10).Value = "Plátano"
$ws.Cells.Item($r, 11).Value = "Sin especificar"
$ws.Cells.Item($r, 12).Value = "Pintón"
$ws.Cells.Item($r, 13).Value = 120
$ws.Cells.Item($r, 14).Value = 15000
$ws.Cells.Item($r, 15).Value = 16000
$ws.Cells.Item($r, 16).Value = 15500
$ws.Cells.Item($r, 17).Value = "$/caja 20 kilos"
$ws.Cells.Item($r, 18).Value = "Bolivia"
$ws.Cells.Item($r, 19).Value = 775
$ws.Cells.Item($r, 20).Value = 20
